$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 193.72223
$ws.Range("I33").Value = 119.07143
$ws.Range("K33").Value = 119.07143
$ws.Range("M33").Value = 109.92857
$ws.Range("H70").Value = 2729.0476
$ws.Range("I70").Value = 811.8
$ws.Range("J70").Value = 3328.1875
$ws.Range("K70").Value = 2435.4
$ws.Range("L70").Value = 9984.5625
$ws.Range("M70").Value = -2165.4
$ws.Range("N70").Value = -10524.5625
$ws.Range("H73").Value = 2729.0476
$ws.Range("I73").Value = 811.8
$ws.Range("J73").Value = 3328.1875
$ws.Range("K73").Value = 2435.4
$ws.Range("L73").Value = 9984.5625
$ws.Range("M73").Value = -1499.4
$ws.Range("N73").Value = -11856.5625
$ws.Range("H74").Value = 5004680.5
$ws.Range("I74").Value = 14289115
$ws.Range("J74").Value = 5369.231
$ws.Range("K74").Value = 14289115
$ws.Range("L74").Value = 5369.231
$ws.Range("M74").Value = -14288179
$ws.Range("N74").Value = -7241.231
$ws.Range("H77").Value = 5004680.5
$ws.Range("I77").Value = 14289115
$ws.Range("J77").Value = 5369.231
$ws.Range("K77").Value = 71445575
$ws.Range("L77").Value = 26846.155
$ws.Range("M77").Value = -71440895
$ws.Range("N77").Value = -36206.155
$ws.Range("H98").Value = 2500.2415
$ws.Range("I98").Value = 1340.28
$ws.Range("J98").Value = 9750
$ws.Range("K98").Value = 1340.28
$ws.Range("L98").Value = 9750
$ws.Range("M98").Value = 157.72
$ws.Range("N98").Value = -12746
$ws.Range("H112").Value = 27028754
$ws.Range("J112").Value = 1813.2354
$ws.Range("L112").Value = 5439.706200000001
$ws.Range("N112").Value = -7655.706200000001
$ws.Range("H122").Value = 2500.2415
$ws.Range("I122").Value = 1340.28
$ws.Range("J122").Value = 9750
$ws.Range("K122").Value = 4020.84
$ws.Range("L122").Value = 29250
$ws.Range("M122").Value = -1570.84
$ws.Range("N122").Value = -34150
$ws.Range("H132").Value = 27784142
$ws.Range("I132").Value = 38468028
$ws.Range("J132").Value = 6039
$ws.Range("K132").Value = 115404084
$ws.Range("L132").Value = 18117
$ws.Range("M132").Value = -115401554
$ws.Range("N132").Value = -23177
$ws.Range("H137").Value = 2512.468
$ws.Range("I137").Value = 1103.5454
$ws.Range("K137").Value = 3310.6362
$ws.Range("M137").Value = -760.6361999999999
$ws.Range("H138").Value = 5141.98
$ws.Range("I138").Value = 745.6667
$ws.Range("J138").Value = 5917.8
$ws.Range("K138").Value = 2237.0001
$ws.Range("L138").Value = 17753.4
$ws.Range("M138").Value = 2902.9999
$ws.Range("N138").Value = -28033.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4689.698
$ws.Range("I32").Value = 4422.827
$ws.Range("J32").Value = 5951.273
$ws.Range("K32").Value = 4422.827
$ws.Range("L32").Value = 5951.273
$ws.Range("M32").Value = -4135.827
$ws.Range("N32").Value = -6525.273
$ws.Range("H61").Value = 1444.7
$ws.Range("I61").Value = 1113.6666
$ws.Range("J61").Value = 2768.8333
$ws.Range("K61").Value = 1113.6666
$ws.Range("L61").Value = 2768.8333
$ws.Range("M61").Value = -901.6666
$ws.Range("N61").Value = -3192.8333
$ws.Range("H74").Value = 3084.973
$ws.Range("I74").Value = 3097.2
$ws.Range("K74").Value = 3097.2
$ws.Range("M74").Value = -2223.2
$ws.Range("H77").Value = 3084.973
$ws.Range("I77").Value = 3097.2
$ws.Range("K77").Value = 15486
$ws.Range("M77").Value = -11118
$ws.Range("H136").Value = 1444.7
$ws.Range("I136").Value = 1113.6666
$ws.Range("J136").Value = 2768.8333
$ws.Range("K136").Value = 3340.9998
$ws.Range("L136").Value = 8306.499899999999
$ws.Range("M136").Value = -790.9998000000001
$ws.Range("N136").Value = -13406.4999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1280.0476
$ws.Range("I86").Value = 1162.5294
$ws.Range("J86").Value = 1779.5
$ws.Range("K86").Value = 1162.5294
$ws.Range("L86").Value = 1779.5
$ws.Range("M86").Value = -39.5293999999999
$ws.Range("N86").Value = -4025.5
$ws.Range("H89").Value = 1280.0476
$ws.Range("I89").Value = 1162.5294
$ws.Range("J89").Value = 1779.5
$ws.Range("K89").Value = 5812.646999999999
$ws.Range("L89").Value = 8897.5
$ws.Range("M89").Value = -196.646999999999
$ws.Range("N89").Value = -20129.5
$ws.Range("H134").Value = 1913.5
$ws.Range("I134").Value = 1208.4529
$ws.Range("J134").Value = 3692.9048
$ws.Range("K134").Value = 3625.3587
$ws.Range("L134").Value = 11078.7144
$ws.Range("M134").Value = -1090.3587
$ws.Range("N134").Value = -16148.7144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26320890
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H34").Value = 26320890
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H58").Value = 1637.6224
$ws.Range("I58").Value = 1545.5916
$ws.Range("J58").Value = 1879.6296
$ws.Range("K58").Value = 1545.5916
$ws.Range("L58").Value = 1879.6296
$ws.Range("M58").Value = -1342.5916
$ws.Range("N58").Value = -2285.6296
$ws.Range("H132").Value = 2448.0732
$ws.Range("I132").Value = 1960.9412
$ws.Range("J132").Value = 4814.143
$ws.Range("K132").Value = 5882.8236
$ws.Range("L132").Value = 14442.429
$ws.Range("M132").Value = -3352.8236
$ws.Range("N132").Value = -19502.429
$ws.Range("H134").Value = 2699.0972
$ws.Range("I134").Value = 2933
$ws.Range("J134").Value = 2309.2593
$ws.Range("K134").Value = 8799
$ws.Range("L134").Value = 6927.777900000001
$ws.Range("M134").Value = -6264
$ws.Range("N134").Value = -11997.7779
$ws.Range("H136").Value = 1637.6224
$ws.Range("I136").Value = 1545.5916
$ws.Range("J136").Value = 1879.6296
$ws.Range("K136").Value = 4636.7748
$ws.Range("L136").Value = 5638.8888
$ws.Range("M136").Value = -2086.7748
$ws.Range("N136").Value = -10738.8888

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3005.2354
$ws.Range("I137").Value = 1043
$ws.Range("J137").Value = 4075.5454
$ws.Range("K137").Value = 3129
$ws.Range("L137").Value = 12226.6362
$ws.Range("M137").Value = 1971
$ws.Range("N137").Value = -22426.6362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10002443
$ws.Range("I80").Value = 14708231
$ws.Range("J80").Value = 2643.125
$ws.Range("K80").Value = 14708231
$ws.Range("L80").Value = 2643.125
$ws.Range("M80").Value = -14707233
$ws.Range("N80").Value = -4639.125
$ws.Range("H83").Value = 10002443
$ws.Range("I83").Value = 14708231
$ws.Range("J83").Value = 2643.125
$ws.Range("K83").Value = 73541155
$ws.Range("L83").Value = 13215.625
$ws.Range("M83").Value = -73536163
$ws.Range("N83").Value = -23199.625
$ws.Range("H100").Value = 39500
$ws.Range("J100").Value = 39500
$ws.Range("L100").Value = 39500
$ws.Range("N100").Value = -41664
$ws.Range("H106").Value = 35000
$ws.Range("J106").Value = 35000
$ws.Range("L106").Value = 35000
$ws.Range("N106").Value = -37524

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1006.9091
$ws.Range("I16").Value = 880.11536
$ws.Range("J16").Value = 1477.8572
$ws.Range("K16").Value = 880.11536
$ws.Range("L16").Value = 1477.8572
$ws.Range("M16").Value = -710.11536
$ws.Range("N16").Value = -1817.8572
$ws.Range("H132").Value = 3930.7666
$ws.Range("I132").Value = 1258.7
$ws.Range("J132").Value = 9274.9
$ws.Range("K132").Value = 3776.1
$ws.Range("L132").Value = 27824.7
$ws.Range("M132").Value = -1246.1
$ws.Range("N132").Value = -32884.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2827.6924
$ws.Range("I122").Value = 1734.5
$ws.Range("K122").Value = 5203.5
$ws.Range("M122").Value = -2753.5
$ws.Range("H132").Value = 6062105.5
$ws.Range("I132").Value = 704.13513
$ws.Range("K132").Value = 2112.40539
$ws.Range("M132").Value = 417.5946100000001
$ws.Range("H136").Value = 2227.074
$ws.Range("I136").Value = 854.7
$ws.Range("J136").Value = 3942.5417
$ws.Range("K136").Value = 2564.1
$ws.Range("L136").Value = 11827.6251
$ws.Range("M136").Value = -14.10000000000036
$ws.Range("N136").Value = -16927.6251
